$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 155, shifting existing rows 155-246 down to 156-247
$ws.Rows(155).Insert()

# Populate the newly inserted row 155 with its data
$newDate = Get-Date -Year 2022 -Month 2 -Day 14 -Hour 0 -Minute 0 -Second 0

$ws.Cells.Item(155,1).Value = 5
$ws.Cells.Item(155,2).Value = 'Macroferia Regional de Talca'
$ws.Cells.Item(155,3).Value = 'Maule'
$ws.Cells.Item(155,4).Value = $newDate
$ws.Cells.Item(155,5).Value = 7
$ws.Cells.Item(155,6).Value = 100112003
$ws.Cells.Item(155,7).Value = 'Ajo'
$ws.Cells.Item(155,8).Value = 'Chino'
$ws.Cells.Item(155,9).Value = 'Primera'
$ws.Cells.Item(155,10).Value = 200
$ws.Cells.Item(155,11).Value = 20000
$ws.Cells.Item(155,12).Value = 20000
$ws.Cells.Item(155,13).Value = 20000
$ws.Cells.Item(155,14).Value = '$/malla 10 kilos'
$ws.Cells.Item(155,15).Value = 'China'
$ws.Cells.Item(155,16).Value = 2000
$ws.Cells.Item(155,17).Value = 10
$ws.Cells.Item(155,18).Value = 'Hortaliza'
